# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.368.90"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.796.66"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.12"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.59"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.796.15"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.28"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.441.28"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.790.81"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.375.61"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.82"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.00"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.71"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.32"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.700"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  +7.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.92"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.92"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.77"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.04"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  -3.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.14"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.745.04"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.996"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.302"
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.68"
$ws.Range("E45").Value = "  +16.08%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.80"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.41"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.93"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "390.08"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.773.54"
$ws.Range("E51").Value = "  +3.80%  "
